# Resort the sheet tabs: move "总计" (the summary/total sheet) so that it
# comes before "2022-Q2" (the detail sheet). This reproduces the tab
# reorder seen in the workbook diff (总计 becomes the first sheet,
# 2022-Q2 becomes the second), without touching any cell data.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2    = $wb.Worksheets.Item("2022-Q2")

# Move "总计" to sit immediately before "2022-Q2".
$wsTotal.Move($wsQ2)
